# Apply the "withElse" template edit:
#  - Move the hidden "_GoBack" bookmark from the m:if field paragraph to the
#    very start of the first paragraph.
#  - Re-split a few text runs so that the single "interesting" word in each
#    sentence is wrapped in <w:proofErr w:type="spellStart"/> / spellEnd
#    markers (mirrors what Word's background spell-checker emits once a
#    document has actually been edited interactively).
#
# We rebuild each affected paragraph's content in one shot via
# Range.InsertXML, which replaces the whole paragraph (text + pPr) with the
# exact WordprocessingML we want - this lets us place proofErr / bookmark
# elements precisely, which the higher level Range/Find API cannot express.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Paragraph 7: "End of demonstration." ---------------------------------
$p7 = $d.Paragraphs(7).Range
$p7.InsertXML('<w:p ' + $wNs + '><w:r><w:t>En</w:t></w:r><w:r><w:t>d</w:t></w:r><w:r><w:t xml:space="preserve"> of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>demonstration</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p>')

# --- Paragraph 5: "The ELSE paragraph." ------------------------------------
$p5 = $d.Paragraphs(5).Range
$p5.InsertXML('<w:p ' + $wNs + '><w:pPr><w:tabs><w:tab w:val="left" w:pos="3119"/></w:tabs></w:pPr><w:r><w:t xml:space="preserve">The </w:t></w:r><w:r><w:t>ELSE</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>paragraph</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p>')

# --- Paragraph 3: "The THEN paragraph." ------------------------------------
$p3 = $d.Paragraphs(3).Range
$p3.InsertXML('<w:p ' + $wNs + '><w:pPr><w:tabs><w:tab w:val="left" w:pos="3119"/></w:tabs></w:pPr><w:r><w:t xml:space="preserve">The THEN </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>paragraph</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p>')

# --- Paragraph 2: the "m:if" field - drop the _GoBack bookmark pair --------
$p2 = $d.Paragraphs(2).Range
$p2.InsertXML('<w:p ' + $wNs + '><w:pPr><w:tabs><w:tab w:val="left" w:pos="3119"/></w:tabs></w:pPr><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve">m:if </w:instrText></w:r><w:r><w:instrText xml:space="preserve">self.name </w:instrText></w:r><w:r><w:instrText>&lt;&gt;</w:instrText></w:r><w:r><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r><w:instrText>' + "'" + '</w:instrText></w:r><w:r><w:instrText>anydsl</w:instrText></w:r><w:r><w:instrText>' + "'" + '</w:instrText></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r></w:p>')

# --- Paragraph 1: "Basic if demonstration :" + the moved _GoBack bookmark -
# NB: the template's original text used a non-breaking space (U+00A0) before
# the colon (French typographic convention); that character is kept as-is,
# it is simply moved into its own trailing run.
$nbsp = [char]0x00A0
$p1 = $d.Paragraphs(1).Range
$p1.InsertXML('<w:p ' + $wNs + '><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">Basic </w:t></w:r><w:r><w:t>if</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>demonstration</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>' + $nbsp + ':</w:t></w:r></w:p>')
